# Applies the "version final runs has added" edit:
#  - Updates D2:D3 (problem_num) and E2:E11 (cpu_time_limit) values
#  - Moves the active cell selection to I19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (problem_num) updates
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1

# Column E (cpu_time_limit) updates
$ws.Range("E2").Value = 100
$ws.Range("E3").Value = 100
$ws.Range("E4").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("E6").Value = 100
$ws.Range("E7").Value = 100
$ws.Range("E8").Value = 100
$ws.Range("E9").Value = 100
$ws.Range("E10").Value = 100
$ws.Range("E11").Value = 100

# Update the active selection cell to match the saved view state
$ws.Range("I19").Select()
